$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 5-7 (the "ECs -> *" rows) entirely; the remaining
# "FAPs -> *" rows (old 5-7) are promoted to rows 2-4 in their place.
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf14"
$ws.Range("C2").Value = "Ltbr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05992133333333333
$ws.Range("H2").Value = 0.179764
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.436245333333333
$ws.Range("N2").Value = 19.308736
$ws.Range("O2").Value = 0.2367562936388591
$ws.Range("P2").Value = 0.2367562936388591
$ws.Range("Q2").Value = 0.3856684020337778
$ws.Range("R2").Value = 3.471015618304
$ws.Range("S2").Value = 0.2367562936388591
$ws.Range("T2").Value = 0.2367562936388591

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf14"
$ws.Range("C3").Value = "Ltbr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05992133333333333
$ws.Range("H3").Value = 0.179764
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.97527
$ws.Range("N3").Value = 44.92581
$ws.Range("O3").Value = 0.5508630013028089
$ws.Range("P3").Value = 0.550863001302809
$ws.Range("Q3").Value = 0.8973381454266667
$ws.Range("R3").Value = 8.076043308839999
$ws.Range("S3").Value = 0.5508630013028089
$ws.Range("T3").Value = 0.550863001302809

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf14"
$ws.Range("C4").Value = "Ltbr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05992133333333333
$ws.Range("H4").Value = 0.179764
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.773592333333333
$ws.Range("N4").Value = 17.320777
$ws.Range("O4").Value = 0.212380705058332
$ws.Range("P4").Value = 0.212380705058332
$ws.Range("Q4").Value = 0.3459613507364445
$ws.Range("R4").Value = 3.113652156628
$ws.Range("S4").Value = 0.212380705058332
$ws.Range("T4").Value = 0.212380705058332
